$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row ---
$ws.Range("E1").Value = "Total-DB"
$ws.Range("F1").Value = "Total-Project"
$ws.Range("G1").Value = "Rate"

# Apply header style (bold/border/centered) from an existing header cell to
# the updated/new header cells E1:G1
$ws.Range("D1").Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 2: Activiti ---
$ws.Range("D2").Value = 200
$ws.Range("E2").Value = 339
$ws.Range("F2").Value = 4890
$ws.Range("G2").Value = 6.932515337423313

# --- Row 3: che ---
$ws.Range("E3").Value = 425
$ws.Range("F3").Value = 2538
$ws.Range("G3").Value = 16.74546887312845

# --- Row 4: pinpoint ---
$ws.Range("D4").Value = 17
$ws.Range("E4").Value = 35
$ws.Range("F4").Value = 8149
$ws.Range("G4").Value = 0.4295005522149957

# --- Row 5: skywalking ---
$ws.Range("E5").Value = 11
$ws.Range("F5").Value = 2587
$ws.Range("G5").Value = 0.4252029377657518

# --- Row 6: storm ---
$ws.Range("E6").Value = 14
$ws.Range("F6").Value = 3398
$ws.Range("G6").Value = 0.4120070629782225

# --- Delete rows 7 through 9 (docker-java, languagetool, ebean) ---
$ws.Range("A7:E9").Delete()
